$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.152962333333333
$ws.Range("H2").Value = 6.458887000000001
$ws.Range("I2").Value = 0.1024890697041326
$ws.Range("J2").Value = 0.1024890697041326
$ws.Range("M2").Value = 9.426699666666666
$ws.Range("N2").Value = 28.280099
$ws.Range("O2").Value = 0.4880118193702016
$ws.Range("P2").Value = 0.4880118193702015
$ws.Range("Q2").Value = 20.29532930997922
$ws.Range("R2").Value = 182.657963789813
$ws.Range("S2").Value = 0.05001587737187314
$ws.Range("T2").Value = 0.05001587737187315

$ws.Range("G3").Value = 2.152962333333333
$ws.Range("H3").Value = 6.458887000000001
$ws.Range("I3").Value = 0.1024890697041326
$ws.Range("J3").Value = 0.1024890697041326
$ws.Range("M3").Value = 7.983522666666666
$ws.Range("O3").Value = 0.4132998355002127
$ws.Range("P3").Value = 0.4132998355002127
$ws.Range("Q3").Value = 17.18822358864622
$ws.Range("R3").Value = 154.694012297816
$ws.Range("S3").Value = 0.04235871564928782
$ws.Range("T3").Value = 0.04235871564928783

$ws.Range("G4").Value = 2.152962333333333
$ws.Range("H4").Value = 6.458887000000001
$ws.Range("I4").Value = 0.1024890697041326
$ws.Range("J4").Value = 0.1024890697041326
$ws.Range("M4").Value = 1.816582333333334
$ws.Range("N4").Value = 5.449747
$ws.Range("O4").Value = 0.09404284435416221
$ws.Range("P4").Value = 0.0940428443541622
$ws.Range("Q4").Value = 3.911033339065445
$ws.Range("R4").Value = 35.19930005158901
$ws.Range("S4").Value = 0.009638363630188621
$ws.Range("T4").Value = 0.009638363630188622

$ws.Range("G5").Value = 2.152962333333333
$ws.Range("H5").Value = 6.458887000000001
$ws.Range("I5").Value = 0.1024890697041326
$ws.Range("J5").Value = 0.1024890697041326
$ws.Range("M5").Value = 0.08973500000000001
$ws.Range("N5").Value = 0.269205
$ws.Range("O5").Value = 0.004645500775423563
$ws.Range("P5").Value = 0.004645500775423562
$ws.Range("Q5").Value = 0.1931960749816667
$ws.Range("R5").Value = 1.738764674835
$ws.Range("S5").Value = 0.0004761130527829875
$ws.Range("T5").Value = 0.0004761130527829875

$ws.Range("I6").Value = 0.0002147722290923241
$ws.Range("J6").Value = 0.0002147722290923242
$ws.Range("M6").Value = 9.426699666666666
$ws.Range("N6").Value = 28.280099
$ws.Range("O6").Value = 0.4880118193702016
$ws.Range("P6").Value = 0.4880118193702015
$ws.Range("Q6").Value = 0.04253012666277777
$ws.Range("R6").Value = 0.382771139965
$ws.Range("S6").Value = 0.0001048113862695389
$ws.Range("T6").Value = 0.0001048113862695389

$ws.Range("I7").Value = 0.0002147722290923241
$ws.Range("J7").Value = 0.0002147722290923242
$ws.Range("M7").Value = 7.983522666666666
$ws.Range("O7").Value = 0.4132998355002127
$ws.Range("P7").Value = 0.4132998355002127
$ws.Range("Q7").Value = 0.03601899309777777
$ws.Range("S7").Value = 0.00008876532695387156
$ws.Range("T7").Value = 0.00008876532695387156

$ws.Range("I8").Value = 0.0002147722290923241
$ws.Range("J8").Value = 0.0002147722290923242
$ws.Range("M8").Value = 1.816582333333334
$ws.Range("N8").Value = 5.449747
$ws.Range("O8").Value = 0.09404284435416221
$ws.Range("P8").Value = 0.0940428443541622
$ws.Range("Q8").Value = 0.008195813960555557
$ws.Range("R8").Value = 0.073762325645
$ws.Range("S8").Value = 0.00002019779131212591
$ws.Range("T8").Value = 0.00002019779131212591

$ws.Range("I9").Value = 0.0002147722290923241
$ws.Range("J9").Value = 0.0002147722290923242
$ws.Range("M9").Value = 0.08973500000000001
$ws.Range("N9").Value = 0.269205
$ws.Range("O9").Value = 0.004645500775423563
$ws.Range("P9").Value = 0.004645500775423562
$ws.Range("Q9").Value = 0.0004048544083333334
$ws.Range("R9").Value = 0.003643689675
$ws.Range("S9").Value = 0.000000997724556787839
$ws.Range("T9").Value = 0.000000997724556787839

$ws.Range("G10").Value = 3.553611
$ws.Range("H10").Value = 10.660833
$ws.Range("I10").Value = 0.1691651915323982
$ws.Range("J10").Value = 0.1691651915323982
$ws.Range("M10").Value = 9.426699666666666
$ws.Range("N10").Value = 28.280099
$ws.Range("O10").Value = 0.4880118193702016
$ws.Range("P10").Value = 0.4880118193702015
$ws.Range("Q10").Value = 33.498823629163
$ws.Range("R10").Value = 301.489412662467
$ws.Range("S10").Value = 0.08255461289383426
$ws.Range("T10").Value = 0.08255461289383426

$ws.Range("G11").Value = 3.553611
$ws.Range("H11").Value = 10.660833
$ws.Range("I11").Value = 0.1691651915323982
$ws.Range("J11").Value = 0.1691651915323982
$ws.Range("M11").Value = 7.983522666666666
$ws.Range("O11").Value = 0.4132998355002127
$ws.Range("P11").Value = 0.4132998355002127
$ws.Range("Q11").Value = 28.370333967016
$ws.Range("R11").Value = 255.333005703144
$ws.Range("S11").Value = 0.06991594583270214
$ws.Range("T11").Value = 0.06991594583270214

$ws.Range("G12").Value = 3.553611
$ws.Range("H12").Value = 10.660833
$ws.Range("I12").Value = 0.1691651915323982
$ws.Range("J12").Value = 0.1691651915323982
$ws.Range("M12").Value = 1.816582333333334
$ws.Range("N12").Value = 5.449747
$ws.Range("O12").Value = 0.09404284435416221
$ws.Range("P12").Value = 0.0940428443541622
$ws.Range("Q12").Value = 6.455426962139001
$ws.Range("R12").Value = 58.09884265925101
$ws.Range("S12").Value = 0.01590877577742336
$ws.Range("T12").Value = 0.01590877577742336

$ws.Range("G13").Value = 3.553611
$ws.Range("H13").Value = 10.660833
$ws.Range("I13").Value = 0.1691651915323982
$ws.Range("J13").Value = 0.1691651915323982
$ws.Range("M13").Value = 0.08973500000000001
$ws.Range("N13").Value = 0.269205
$ws.Range("O13").Value = 0.004645500775423563
$ws.Range("P13").Value = 0.004645500775423562
$ws.Range("Q13").Value = 0.318883283085
$ws.Range("R13").Value = 2.869949547765001
$ws.Range("S13").Value = 0.0007858570284384313
$ws.Range("T13").Value = 0.0007858570284384313

$ws.Range("G14").Value = 15.29566566666667
$ws.Range("H14").Value = 45.886997
$ws.Range("I14").Value = 0.7281309665343768
$ws.Range("J14").Value = 0.7281309665343769
$ws.Range("M14").Value = 9.426699666666666
$ws.Range("N14").Value = 28.280099
$ws.Range("O14").Value = 0.4880118193702016
$ws.Range("P14").Value = 0.4880118193702015
$ws.Range("Q14").Value = 144.1876464414114
$ws.Range("R14").Value = 1297.688817972703
$ws.Range("S14").Value = 0.3553365177182246
$ws.Range("T14").Value = 0.3553365177182246

$ws.Range("G15").Value = 15.29566566666667
$ws.Range("H15").Value = 45.886997
$ws.Range("I15").Value = 0.7281309665343768
$ws.Range("J15").Value = 0.7281309665343769
$ws.Range("M15").Value = 7.983522666666666
$ws.Range("O15").Value = 0.4132998355002127
$ws.Range("P15").Value = 0.4132998355002127
$ws.Range("Q15").Value = 122.1132935515884
$ws.Range("R15").Value = 1099.019641964296
$ws.Range("S15").Value = 0.3009364086912688
$ws.Range("T15").Value = 0.3009364086912688

$ws.Range("G16").Value = 15.29566566666667
$ws.Range("H16").Value = 45.886997
$ws.Range("I16").Value = 0.7281309665343768
$ws.Range("J16").Value = 0.7281309665343769
$ws.Range("M16").Value = 1.816582333333334
$ws.Range("N16").Value = 5.449747
$ws.Range("O16").Value = 0.09404284435416221
$ws.Range("P16").Value = 0.0940428443541622
$ws.Range("Q16").Value = 27.78583602663989
$ws.Range("R16").Value = 250.072524239759
$ws.Range("S16").Value = 0.06847550715523809
$ws.Range("T16").Value = 0.06847550715523809

$ws.Range("G17").Value = 15.29566566666667
$ws.Range("H17").Value = 45.886997
$ws.Range("I17").Value = 0.7281309665343768
$ws.Range("J17").Value = 0.7281309665343769
$ws.Range("M17").Value = 0.08973500000000001
$ws.Range("N17").Value = 0.269205
$ws.Range("O17").Value = 0.004645500775423563
$ws.Range("P17").Value = 0.004645500775423562
$ws.Range("Q17").Value = 1.372556558598333
$ws.Range("R17").Value = 12.353009027385
$ws.Range("S17").Value = 0.003382532969645356
$ws.Range("T17").Value = 0.003382532969645356

